# Add a new "amazonTest" test case row to the RUNMANAGER sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RUNMANAGER")

$ws.Range("A4").Value = "amazonTest"
$ws.Range("B4").Value = "To verify if the amazon test is working or not"
$ws.Range("C4").Value = "yes"
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 1

# Reflect the final selection/active-sheet state left by the editor.
$dataWs = $wb.Worksheets.Item("DATA")
$dataWs.Range("A8").Select() | Out-Null

$ws.Activate() | Out-Null
$ws.Range("F14").Select() | Out-Null
